# Update cryptocurrency price (column D) and volume/change % (column E) cells
# to reflect the latest scraped values, per the GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '62.688.86'
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.63%  '
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.453.65'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.76%  '
$ws.Cells.Item(4, 5).Value = '  +0.07%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '570.63'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.21%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '145.81'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.68%  '
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.527'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -2.04%  '
$ws.Cells.Item(9, 5).Value = '  -1.35%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.15'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -2.31%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.347'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -1.50%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '28.58'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -1.09%  '
$ws.Cells.Item(14, 5).Value = '  -3.20%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.895.59'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -0.86%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '62.556.24'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.64%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.453.76'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.72%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.65'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -6.60%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.70'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -2.93%  '
$ws.Cells.Item(20, 5).Value = '  -0.26%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '321.01'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -2.69%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.12'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.33%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.03%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.86'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +2.77%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '64.67'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -2.40%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '647.36'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -2.39%  '
$ws.Cells.Item(27, 5).Value = '  -0.54%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0948'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -4.17%  '
$ws.Cells.Item(29, 5).Value = '  -0.07%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.40'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -3.07%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.80'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -3.26%  '
$ws.Cells.Item(32, 5).Value = '  -3.41%  '
$ws.Cells.Item(33, 5).Value = '  +0.05%  '
$ws.Cells.Item(34, 5).Value = '  -0.01%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.47'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -4.40%  '
$ws.Cells.Item(36, 5).Value = '  -3.35%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '150.48'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -1.30%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.52'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -1.43%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.362'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -2.46%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.60'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -4.48%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.68'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -4.10%  '
$ws.Cells.Item(43, 5).Value = '  -0.01%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0₆0305'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -0.59%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '152.35'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.29%  '
$ws.Cells.Item(46, 5).Value = '  +1.76%  '
$ws.Cells.Item(47, 5).Value = '  -2.30%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.601'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.76%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.87'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -3.62%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0503'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.80%  '
$ws.Cells.Item(51, 5).Value = '  -2.12%  '
